# Append three new parking records (rows 57-59) to the Data_Parking sheet,
# mirroring the data added in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A="QIPZ5P2BG6MW"; B="D 5034 YAW"; C="Motor"; D="2025-01-13 17:00:18"; E="2025-01-13 17:00:57"; F="0.0 jam 0.0 menit 39.717722 detik"; G=2000; H="Reza Ramdan Permana"; I="./capture/masuk/QIPZ5P2BG6MW.png"; J="./capture/keluar/QIPZ5P2BG6MW.png" },
    @{ A="8E3FKYN6OBLV"; B="D 5430 BAH"; C="Mobil"; D="2025-01-13 17:14:36"; E="2025-01-13 17:15:17"; F="00:00:41"; G=4000; H="Reza Ramdan Permana"; I="./capture/masuk/8E3FKYN6OBLV.png"; J="./capture/keluar/8E3FKYN6OBLV.png" },
    @{ A="3GH15R2GL0YE"; B="D 6405 HWK"; C="Motor"; D="2025-01-13 17:21:11"; E="2025-01-13 17:22:06"; F="00:00:55"; G=2000; H="Reza Ramdan Permana"; I="./capture/masuk/3GH15R2GL0YE.png"; J="./capture/keluar/3GH15R2GL0YE.png" }
)

$startRow = 57
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
}
